$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------------
# 1) Insert a new row at row 24. This pushes the old "12filters/3b(0.9)" block
#    (rows 24-30) down to rows 25-31, and the old "3b(v2)" block (rows 32-34)
#    down to rows 33-35 -- matching the target layout. The existing merged
#    cells (A24:A30 -> A25:A31, A32:A34 -> A33:A35) are kept in sync
#    automatically by Excel.
# ---------------------------------------------------------------------------
$ws.Rows("24:24").Insert()

# ---------------------------------------------------------------------------
# 2) Fill in the newly inserted row 24 with the "20% (excl input)" checkpoint,
#    highlighted in yellow like the other 0.2-noise-level row (row 20).
# ---------------------------------------------------------------------------
$ws.Range("A17").Copy()
$ws.Range("A24").PasteSpecial($xlPasteFormats)
$ws.Range("B16").Copy()
$ws.Range("B24").PasteSpecial($xlPasteFormats)
$ws.Range("C20:F20").Copy()
$ws.Range("C24:F24").PasteSpecial($xlPasteFormats)

$ws.Range("B24").Value = "20% (excl input)"
$ws.Range("C24").Value = 0.85750000000000004
$ws.Range("D24").Value = 0.85399999999999998
$ws.Range("E24").Value = 0.83779999999999999
$ws.Range("F24").Value = 0.80100000000000005

# ---------------------------------------------------------------------------
# 3) Add the new "16filter_3b" checkpoint group at rows 37-40.
# ---------------------------------------------------------------------------
$ws.Range("A25").Copy()
$ws.Range("A37").PasteSpecial($xlPasteFormats)
$ws.Range("A26").Copy()
$ws.Range("A38:A40").PasteSpecial($xlPasteFormats)
$ws.Range("B26:F26").Copy()
$ws.Range("B37:F37").PasteSpecial($xlPasteFormats)
$ws.Range("B37:F37").Copy()
$ws.Range("B38:F40").PasteSpecial($xlPasteFormats)

$ws.Range("A37").Value = "16filter_3b"
$ws.Range("B37").Value = 0.15
$ws.Range("C37").Value = 0.86529999999999996
$ws.Range("D37").Value = 0.86329999999999996
$ws.Range("E37").Value = 0.84089999999999998
$ws.Range("F37").Value = 0.79849999999999999

$ws.Range("B38").Value = 0.2
$ws.Range("C38").Value = 0.86380000000000001
$ws.Range("D38").Value = 0.8579
$ws.Range("E38").Value = 0.84750000000000003
$ws.Range("F38").Value = 0.82320000000000004

$ws.Range("B39").Value = "0.15 (excl input)"
$ws.Range("C39").Value = 0.88219999999999998
$ws.Range("D39").Value = 0.87680000000000002
$ws.Range("E39").Value = 0.86150000000000004
$ws.Range("F39").Value = 0.82499999999999996

$ws.Range("B40").Value = "0.2 (excl input)"
$ws.Range("C40").Value = 0.87819999999999998
$ws.Range("D40").Value = 0.87549999999999994
$ws.Range("E40").Value = 0.86250000000000004
$ws.Range("F40").Value = 0.83560000000000001

$ws.Range("A37:A40").Merge()

# ---------------------------------------------------------------------------
# 4) Fill in row 23 (previously an empty/gap row) with the new
#    "15%(excl input)" checkpoint belonging to the 16-filter (3b 0.95) group.
# ---------------------------------------------------------------------------
$ws.Range("A17").Copy()
$ws.Range("A23").PasteSpecial($xlPasteFormats)
$ws.Range("B16").Copy()
$ws.Range("B23").PasteSpecial($xlPasteFormats)
$ws.Range("C17:F17").Copy()
$ws.Range("C23:F23").PasteSpecial($xlPasteFormats)

$ws.Range("B23").Value = "15%(excl input)"
$ws.Range("C23").Value = 0.86119999999999997
$ws.Range("D23").Value = 0.85399999999999998
$ws.Range("E23").Value = 0.82950000000000002
$ws.Range("F23").Value = 0.78320000000000001

# ---------------------------------------------------------------------------
# 5) Re-merge the "3b (0.95)" group label cell now that it spans two more
#    rows (A16:A22 -> A16:A24).
# ---------------------------------------------------------------------------
$ws.Range("A16:A22").UnMerge()
$ws.Range("A16:A24").Merge()

# ---------------------------------------------------------------------------
# 6) Restore the saved cursor/selection position.
# ---------------------------------------------------------------------------
$ws.Range("F24").Select()
